# GDP q3 update haverpull
# Revises a batch of previously-published monthly figures (small
# haver-pull corrections scattered through the history) and appends
# the new September 2023 row (row 646) with an adjusted August 2023
# row (row 645: new lwcl/lufp figures plus other corrected columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F98").Value = 13886979
$ws.Range("G98").Value = 11076413
$ws.Range("H98").Value = 904873
$ws.Range("M98").Value = 8643095
$ws.Range("M99").Value = 8586182
$ws.Range("M101").Value = 8496279
$ws.Range("M102").Value = 8469982
$ws.Range("M103").Value = 8390320
$ws.Range("M104").Value = 8353745
$ws.Range("M105").Value = 8358622
$ws.Range("M106").Value = 8313965
$ws.Range("M107").Value = 8307260
$ws.Range("M108").Value = 8263909
$ws.Range("M109").Value = 8212244
$ws.Range("H132").Value = 935419
$ws.Range("M132").Value = 14066966
$ws.Range("M133").Value = 14469419
$ws.Range("M134").Value = 14605054
$ws.Range("M135").Value = 14690491
$ws.Range("M136").Value = 14866950
$ws.Range("M137").Value = 14861468
$ws.Range("M138").Value = 14670860
$ws.Range("M139").Value = 14467085
$ws.Range("D140").Value = 823677
$ws.Range("M140").Value = 14131865
$ws.Range("N140").Value = 8935513
$ws.Range("M141").Value = 13894161
$ws.Range("N141").Value = 8782028
$ws.Range("M142").Value = 13750068
$ws.Range("N142").Value = 8776296
$ws.Range("M143").Value = 13621937
$ws.Range("N143").Value = 8836738
$ws.Range("N144").Value = 9042546
$ws.Range("N145").Value = 9395166
$ws.Range("N146").Value = 9716030
$ws.Range("N147").Value = 9916861
$ws.Range("N148").Value = 10176494
$ws.Range("N149").Value = 10402075
$ws.Range("N150").Value = 10586475
$ws.Range("N151").Value = 10834388
$ws.Range("D186").Value = 590421
$ws.Range("G186").Value = 9817862
$ws.Range("H186").Value = 1209128
$ws.Range("I186").Value = 127.38
$ws.Range("M186").Value = 13889917
$ws.Range("N186").Value = 8211226
$ws.Range("M187").Value = 13966150
$ws.Range("N187").Value = 8231217
$ws.Range("M188").Value = 14177384
$ws.Range("N188").Value = 8311003
$ws.Range("M189").Value = 14260691
$ws.Range("N189").Value = 8315363
$ws.Range("M190").Value = 14423018
$ws.Range("N190").Value = 8370944
$ws.Range("D191").Value = 573555
$ws.Range("G191").Value = 8670559
$ws.Range("H191").Value = 1073078
$ws.Range("M191").Value = 14528194
$ws.Range("N191").Value = 8368473
$ws.Range("M192").Value = 14532968
$ws.Range("N192").Value = 8290047
$ws.Range("M193").Value = 14761492
$ws.Range("N193").Value = 8365139
$ws.Range("F194").Value = 14729048
$ws.Range("M194").Value = 14844173
$ws.Range("N194").Value = 8388026
$ws.Range("M195").Value = 14889103
$ws.Range("N195").Value = 8325498
$ws.Range("F196").Value = 13097754
$ws.Range("M196").Value = 14967785
$ws.Range("N196").Value = 8319994
$ws.Range("M197").Value = 15089929
$ws.Range("N197").Value = 8341458
$ws.Range("M198").Value = 15141526
$ws.Range("N198").Value = 8330315
$ws.Range("M199").Value = 15309070
$ws.Range("N199").Value = 8356563
$ws.Range("M200").Value = 15440680
$ws.Range("N200").Value = 8350582
$ws.Range("M201").Value = 15510454
$ws.Range("N201").Value = 8368325
$ws.Range("F202").Value = 10071179
$ws.Range("M202").Value = 15673764
$ws.Range("N202").Value = 8381882
$ws.Range("F203").Value = 9851360
$ws.Range("B232").Value = 1361139
$ws.Range("F232").Value = 11589036
$ws.Range("B636").Value = 1008284
$ws.Range("D636").Value = 353166
$ws.Range("E636").Value = 112106
$ws.Range("F636").Value = 5706299
$ws.Range("G636").Value = 4582761
$ws.Range("H636").Value = 1871519
$ws.Range("I636").Value = 417.45
$ws.Range("M636").Value = 24210263
$ws.Range("N636").Value = 4314198
$ws.Range("O636").Value = 1586005
$ws.Range("B637").Value = 1159588
$ws.Range("D637").Value = 463463
$ws.Range("E637").Value = 114578
$ws.Range("F637").Value = 6575471
$ws.Range("G637").Value = 5403279
$ws.Range("H637").Value = 2237864
$ws.Range("I637").Value = 422.41
$ws.Range("J637").Value = 31.21
$ws.Range("M637").Value = 24255635
$ws.Range("N637").Value = 4362818
$ws.Range("O637").Value = 1518550
$ws.Range("B638").Value = 1226023
$ws.Range("D638").Value = 650582
$ws.Range("E638").Value = 140525
$ws.Range("F638").Value = 8985035
$ws.Range("G638").Value = 7350768
$ws.Range("H638").Value = 3117353
$ws.Range("I638").Value = 430.24
$ws.Range("J638").Value = 31.91
$ws.Range("M638").Value = 24737814
$ws.Range("N638").Value = 4483545
$ws.Range("O638").Value = 1481373
$ws.Range("B639").Value = 898590
$ws.Range("D639").Value = 467248
$ws.Range("E639").Value = 118662
$ws.Range("F639").Value = 7771291
$ws.Range("G639").Value = 6556427
$ws.Range("H639").Value = 2820505
$ws.Range("I639").Value = 437.6
$ws.Range("J639").Value = 32.19
$ws.Range("M639").Value = 25142790
$ws.Range("N639").Value = 4558158
$ws.Range("O639").Value = 1460858
$ws.Range("B640").Value = 986573
$ws.Range("D640").Value = 410593
$ws.Range("E640").Value = 130461
$ws.Range("F640").Value = 7972699
$ws.Range("G640").Value = 6869973
$ws.Range("H640").Value = 2977061
$ws.Range("I640").Value = 441.08
$ws.Range("L640").Value = 413.66
$ws.Range("M640").Value = 25579137
$ws.Range("N640").Value = 4598664
$ws.Range("O640").Value = 1422862
$ws.Range("J641").Value = 33.71
$ws.Range("L641").Value = 416.42
$ws.Range("M641").Value = 26306935
$ws.Range("N641").Value = 4657554
$ws.Range("O641").Value = 1446224
$ws.Range("B642").Value = 1002759
$ws.Range("D642").Value = 386133
$ws.Range("E642").Value = 155727
$ws.Range("F642").Value = 7189797
$ws.Range("G642").Value = 6150605
$ws.Range("H642").Value = 2648685
$ws.Range("I642").Value = 438.29
$ws.Range("J642").Value = 34.22
$ws.Range("L642").Value = 419.82
$ws.Range("M642").Value = 27140503
$ws.Range("N642").Value = 4730666
$ws.Range("O642").Value = 1476376
$ws.Range("B643").Value = 1113019
$ws.Range("D643").Value = 417214
$ws.Range("E643").Value = 135738
$ws.Range("F643").Value = 6939193
$ws.Range("G643").Value = 5719857
$ws.Range("H643").Value = 2439855
$ws.Range("I643").Value = 434.31
$ws.Range("J643").Value = 34.51
$ws.Range("L643").Value = 422.68
$ws.Range("M643").Value = 27887979
$ws.Range("N643").Value = 4820005
$ws.Range("O643").Value = 1505789
$ws.Range("R643").Value = 396538
$ws.Range("B644").Value = 1058702
$ws.Range("D644").Value = 467105
$ws.Range("E644").Value = 165492
$ws.Range("F644").Value = 8488625
$ws.Range("G644").Value = 6749139
$ws.Range("H644").Value = 2807083
$ws.Range("I644").Value = 423.46
$ws.Range("J644").Value = 34.74
$ws.Range("L644").Value = 425.35
$ws.Range("M644").Value = 28941187
$ws.Range("N644").Value = 4919920
$ws.Range("O644").Value = 1557807
$ws.Range("P644").Value = 5243
$ws.Range("Q644").Value = 14561
$ws.Range("R644").Value = 397510
$ws.Range("B645").Value = 938609
$ws.Range("C645").Value = 1793.4
$ws.Range("D645").Value = 422943
$ws.Range("E645").Value = 156311
$ws.Range("F645").Value = 7609986
$ws.Range("G645").Value = 6595457
$ws.Range("H645").Value = 2808780
$ws.Range("I645").Value = 433.61
$ws.Range("J645").Value = 35
$ws.Range("K645").Value = 14.17
$ws.Range("L645").Value = 428.68
$ws.Range("M645").Value = 29733264
$ws.Range("N645").Value = 4986272
$ws.Range("O645").Value = 1595385
$ws.Range("P645").Value = 5267
$ws.Range("Q645").Value = 14576
$ws.Range("R645").Value = 397601
$ws.Range("A646").Value = 45199
$ws.Range("B646").Value = 739107
$ws.Range("D646").Value = 287345
$ws.Range("E646").Value = 200527
$ws.Range("F646").Value = 6259636
$ws.Range("G646").Value = 5143477
$ws.Range("H646").Value = 2304278
$ws.Range("I646").Value = 456.77
$ws.Range("J646").Value = 36.82
$ws.Range("K646").Value = 14.36
$ws.Range("L646").Value = 432.12
$ws.Range("M646").Value = 30380213
$ws.Range("N646").Value = 4985522
$ws.Range("O646").Value = 1693096
$ws.Range("P646").Value = 5296
$ws.Range("Q646").Value = 14614
